# "fixed some missing sites"
# Row 9 was missing its site_name / type_of_fish / measurement_type /
# measurement_units values (columns A-D). Every other data row (2-8)
# already carries the same values for these columns, so fill row 9 in to
# match them, copying both the values and the cell formatting from row 8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value2 = $ws.Range("A8").Value2
$ws.Range("B9").Value2 = $ws.Range("B8").Value2
$ws.Range("C9").Value2 = $ws.Range("C8").Value2
$ws.Range("D9").Value2 = $ws.Range("D8").Value2

$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Reflect the saved sheet view's active selection.
$ws.Range("D15").Select()
